$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "NOUBAIL MOHAMMED"
$ws.Range("B2").Value = "IR801997"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "007400000313200019604463"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E2").Value = "AWB"
$ws.Range("G2").Value = "002/TTT"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2000

# Update row 3
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2000

# Delete rows 4 and 5
$ws.Range("A4:K5").EntireRow.Delete()
